# LFM Test Plan update: mark quantization-tier retests as Passed/Failed
# with a Last Run Date, and record the new selection left behind by the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tiers1_3")

# Rows whose Status (col F) moves from "Planned" to "Passed", and gets a
# Last Run Date (col G) of 10/23/2025 (serial 45953). Row 28 failed instead.
$passedRows = @(7, 8, 9, 14, 23, 24, 25, 26, 27, 29, 30)
$failedRows = @(28)

foreach ($r in $passedRows) {
    $ws.Range("F$r").Value = "Passed"
    $ws.Range("G$r").Value = 45953
}

foreach ($r in $failedRows) {
    $ws.Range("F$r").Value = "Failed"
    $ws.Range("G$r").Value = 45953
}

# Copy the existing date formatting (style s="5", built-in m/d/yyyy format)
# from G2 onto the newly populated Last Run Date cells instead of letting a
# fresh numFmt get created. Paste one cell at a time - a multi-area union
# range only receives the format on its first area.
$ws.Range("G2").Copy() | Out-Null
foreach ($r in ($passedRows + $failedRows)) {
    $ws.Range("G$r").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
}
$excel.CutCopyMode = $false

# Update the active selection left on the sheet to G23:G30 (anchor G23),
# matching where the author ended up after filling in the new results.
$ws.Activate() | Out-Null
$ws.Range("G23:G30").Select() | Out-Null

Write-Output "Updated status/date for rows $($passedRows -join ', ') (Passed) and $($failedRows -join ', ') (Failed)"
